# Week 9 Meeting.pptx - "added timelog file to repo"
#
# The deck was reopened/re-saved a couple weeks later (11/29/2021 -> 12/11/2021).
# The only user-visible effect baked into the canonical OOXML is that every
# auto-updating "Date Placeholder" (master + all slide layouts) recalculated
# its cached date text to the new save date. Update each of them here.

$p = $ppt.ActivePresentation
$newDate = "12/11/2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try {
            $phType = $sh.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        if ($phType -eq $ppPlaceholderDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}
